$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename TimeSlice labels in column D (rows 13-16): 1R1MOH0N -> 1R1MO{N}H
$ws.Range("D13").Value = "1R1MO1H"
$ws.Range("D14").Value = "1R1MO2H"
$ws.Range("D15").Value = "1R1MO3H"
$ws.Range("D16").Value = "1R1MO4H"

# Update the active selection to G4
$ws.Range("G4").Select()
